# Auto-generated edit script: refresh market-price columns (H:N)
# across all 8 item sheets, per the scheduled-runner commit.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 99
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 0
# Row 116
$ws.Range("H116").Value = 4129.625
$ws.Range("I116").Value = 3770.3635
$ws.Range("J116").Value = 4920
$ws.Range("K116").Value = 3770.3635
$ws.Range("L116").Value = 4920
$ws.Range("M116").Value = -328.3634999999999
$ws.Range("N116").Value = -11804
# Row 132
$ws.Range("H132").Value = 931.57776
$ws.Range("I132").Value = 930.0227
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 2790.0681
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -260.0681
$ws.Range("N132").Value = -8060
# Row 133
$ws.Range("H133").Value = 100000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 100000
$ws.Range("K133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("M133").Value = 100000
$ws.Range("N133").Value = -110120
# Row 135
$ws.Range("H135").Value = 623.5
$ws.Range("I135").Value = 637.6486
$ws.Range("J135").Value = 100
$ws.Range("K135").Value = 5738.8374
$ws.Range("L135").Value = 900
$ws.Range("M135").Value = -3203.8374
$ws.Range("N135").Value = -5970
# Row 137
$ws.Range("H137").Value = 2268.5283
$ws.Range("I137").Value = 2335.475
$ws.Range("J137").Value = 2062.5386
$ws.Range("K137").Value = 7006.424999999999
$ws.Range("L137").Value = 6187.6158
$ws.Range("M137").Value = -4456.424999999999
$ws.Range("N137").Value = -11287.6158

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 14706950
$ws.Range("I2").Value = 16667518
$ws.Range("J2").Value = 2687.5
$ws.Range("K2").Value = 16667518
$ws.Range("L2").Value = 2687.5
$ws.Range("M2").Value = -16667405
$ws.Range("N2").Value = -2913.5
# Row 32
$ws.Range("H32").Value = 9838.096
$ws.Range("I32").Value = 9791.34
$ws.Range("J32").Value = 11007
$ws.Range("K32").Value = 9791.34
$ws.Range("L32").Value = 11007
$ws.Range("M32").Value = -9504.34
$ws.Range("N32").Value = -11581
# Row 45
$ws.Range("H45").Value = 3642.353
$ws.Range("I45").Value = 2842.8
$ws.Range("J45").Value = 4784.5713
$ws.Range("K45").Value = 2842.8
$ws.Range("L45").Value = 4784.5713
$ws.Range("M45").Value = -2465.8
$ws.Range("N45").Value = -5538.5713
# Row 61
$ws.Range("H61").Value = 5635.75
$ws.Range("I61").Value = 7010.6294
$ws.Range("J61").Value = 1511.1111
$ws.Range("K61").Value = 7010.6294
$ws.Range("L61").Value = 1511.1111
$ws.Range("M61").Value = -6798.6294
$ws.Range("N61").Value = -1935.1111
# Row 74
$ws.Range("H74").Value = 3456.9
$ws.Range("I74").Value = 2174.3333
$ws.Range("J74").Value = 15000
$ws.Range("K74").Value = 2174.3333
$ws.Range("L74").Value = 15000
$ws.Range("M74").Value = -1300.3333
$ws.Range("N74").Value = -16748
# Row 77
$ws.Range("H77").Value = 3456.9
$ws.Range("I77").Value = 2174.3333
$ws.Range("J77").Value = 15000
$ws.Range("K77").Value = 10871.6665
$ws.Range("L77").Value = 75000
$ws.Range("M77").Value = -6503.666499999999
$ws.Range("N77").Value = -83736
# Row 110
$ws.Range("H110").Value = 5042
$ws.Range("I110").Value = 2867.3
$ws.Range("J110").Value = 8666.5
$ws.Range("K110").Value = 2867.3
$ws.Range("L110").Value = 8666.5
$ws.Range("M110").Value = -822.3000000000002
$ws.Range("N110").Value = -12756.5
# Row 116
$ws.Range("H116").Value = 14706950
$ws.Range("I116").Value = 16667518
$ws.Range("J116").Value = 2687.5
$ws.Range("K116").Value = 16667518
$ws.Range("L116").Value = 2687.5
$ws.Range("M116").Value = -16665224
$ws.Range("N116").Value = -7275.5
# Row 132
$ws.Range("H132").Value = 1502.3684
$ws.Range("I132").Value = 1502.5
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 4507.5
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -1977.5
$ws.Range("N132").Value = -9560
# Row 136
$ws.Range("H136").Value = 5635.75
$ws.Range("I136").Value = 7010.6294
$ws.Range("J136").Value = 1511.1111
$ws.Range("K136").Value = 21031.8882
$ws.Range("L136").Value = 4533.3333
$ws.Range("M136").Value = -18481.8882
$ws.Range("N136").Value = -9633.3333

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 14706950
$ws.Range("I3").Value = 16667518
$ws.Range("J3").Value = 2687.5
$ws.Range("K3").Value = 16667518
$ws.Range("L3").Value = 2687.5
$ws.Range("M3").Value = -16667404
$ws.Range("N3").Value = -2915.5
# Row 107
$ws.Range("H107").Value = 4742.75
$ws.Range("I107").Value = 4742.75
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 4742.75
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -2822.75
# Row 134
$ws.Range("H134").Value = 3578.818
$ws.Range("I134").Value = 3263.111
$ws.Range("J134").Value = 4999.5
$ws.Range("K134").Value = 9789.332999999999
$ws.Range("L134").Value = 14998.5
$ws.Range("M134").Value = -7254.332999999999
$ws.Range("N134").Value = -20068.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 19998
$ws.Range("I6").Value = 19998
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 19998
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -19885
$ws.Range("N6").ClearContents()
# Row 31
$ws.Range("H31").Value = 29415882
$ws.Range("I31").Value = 50002550
$ws.Range("J31").Value = 6349.9287
$ws.Range("K31").Value = 50002550
$ws.Range("L31").Value = 6349.9287
$ws.Range("M31").Value = -50002255
$ws.Range("N31").Value = -6939.9287
# Row 34
$ws.Range("H34").Value = 29415882
$ws.Range("I34").Value = 50002550
$ws.Range("J34").Value = 6349.9287
$ws.Range("K34").Value = 50002550
$ws.Range("L34").Value = 6349.9287
$ws.Range("M34").Value = -50002348
$ws.Range("N34").Value = -6753.9287
# Row 58
$ws.Range("H58").Value = 12017.667
$ws.Range("I58").Value = 5895.6665
$ws.Range("J58").Value = 14466.467
$ws.Range("K58").Value = 5895.6665
$ws.Range("L58").Value = 14466.467
$ws.Range("M58").Value = -5692.6665
$ws.Range("N58").Value = -14872.467
# Row 99
$ws.Range("H99").Value = 6168.5415
$ws.Range("I99").Value = 6941.5
$ws.Range("J99").Value = 5395.5835
$ws.Range("K99").Value = 6941.5
$ws.Range("L99").Value = 5395.5835
$ws.Range("M99").Value = -5443.5
$ws.Range("N99").Value = -8391.583500000001
# Row 126
$ws.Range("H126").Value = 6168.5415
$ws.Range("I126").Value = 6941.5
$ws.Range("J126").Value = 5395.5835
$ws.Range("K126").Value = 20824.5
$ws.Range("L126").Value = 16186.7505
$ws.Range("M126").Value = -18354.5
$ws.Range("N126").Value = -21126.7505
# Row 132
$ws.Range("H132").Value = 4294.091
$ws.Range("I132").Value = 3344
$ws.Range("J132").Value = 5666.4443
$ws.Range("K132").Value = 10032
$ws.Range("L132").Value = 16999.3329
$ws.Range("M132").Value = -7502
$ws.Range("N132").Value = -22059.3329
# Row 134
$ws.Range("H134").Value = 2777.8286
$ws.Range("I134").Value = 1688.24
$ws.Range("J134").Value = 5501.8
$ws.Range("K134").Value = 5064.72
$ws.Range("L134").Value = 16505.4
$ws.Range("M134").Value = -2529.72
$ws.Range("N134").Value = -21575.4
# Row 136
$ws.Range("H136").Value = 12017.667
$ws.Range("I136").Value = 5895.6665
$ws.Range("J136").Value = 14466.467
$ws.Range("K136").Value = 17686.9995
$ws.Range("L136").Value = 43399.401
$ws.Range("M136").Value = -15136.9995
$ws.Range("N136").Value = -48499.401

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 62.666668
$ws.Range("I2").Value = 66.89474
$ws.Range("J2").Value = 22.5
$ws.Range("K2").Value = 401.36844
$ws.Range("L2").Value = 135
$ws.Range("M2").Value = -288.36844
$ws.Range("N2").Value = -361
# Row 34
$ws.Range("H34").Value = 67192.60000000001
$ws.Range("I34").Value = 573.3333
$ws.Range("J34").Value = 111605.445
$ws.Range("K34").Value = 1719.9999
$ws.Range("L34").Value = 334816.335
$ws.Range("M34").Value = -1635.9999
$ws.Range("N34").Value = -334984.335
# Row 39
$ws.Range("H39").Value = 2027
$ws.Range("I39").Value = 3094.5
$ws.Range("J39").Value = 1600
$ws.Range("K39").Value = 9283.5
$ws.Range("L39").Value = 4800
$ws.Range("M39").Value = -8989.5
$ws.Range("N39").Value = -5388
# Row 70
$ws.Range("H70").Value = 3525
$ws.Range("I70").Value = 1250
$ws.Range("J70").Value = 5800
$ws.Range("K70").Value = 3750
$ws.Range("L70").Value = 17400
$ws.Range("M70").Value = -3435
$ws.Range("N70").Value = -18030
# Row 73
$ws.Range("H73").Value = 3525
$ws.Range("I73").Value = 1250
$ws.Range("J73").Value = 5800
$ws.Range("K73").Value = 3750
$ws.Range("L73").Value = 17400
$ws.Range("M73").Value = -2658
$ws.Range("N73").Value = -19584
# Row 117
$ws.Range("H117").Value = 1819170.6
$ws.Range("I117").Value = 3334161.8
$ws.Range("J117").Value = 1251049
$ws.Range("K117").Value = 10002485.4
$ws.Range("L117").Value = 3753147
$ws.Range("M117").Value = -9999043.399999999
$ws.Range("N117").Value = -3760031
# Row 138
$ws.Range("H138").Value = 4195.4546
$ws.Range("I138").Value = 3905.7778
$ws.Range("J138").Value = 5499
$ws.Range("K138").Value = 11717.3334
$ws.Range("L138").Value = 16497
$ws.Range("M138").Value = -6577.3334
$ws.Range("N138").Value = -26777

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 2626.9375
$ws.Range("I113").Value = 2310.077
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 2310.077
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -140.0770000000002
$ws.Range("N113").Value = -8340
# Row 124
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").ClearContents()
$ws.Range("N124").Value = 0
# Row 126
$ws.Range("H126").Value = 4453.1177
$ws.Range("I126").Value = 4012.875
$ws.Range("J126").Value = 4844.4443
$ws.Range("K126").Value = 12038.625
$ws.Range("L126").Value = 14533.3329
$ws.Range("M126").Value = -9568.625
$ws.Range("N126").Value = -19473.3329
# Row 132
$ws.Range("H132").Value = 1838.8
$ws.Range("I132").Value = 1049.8695
$ws.Range("J132").Value = 3350.9167
$ws.Range("K132").Value = 3149.6085
$ws.Range("L132").Value = 10052.7501
$ws.Range("M132").Value = -619.6085000000003
$ws.Range("N132").Value = -15112.7501

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1199.7142
$ws.Range("I22").Value = 1699.5
$ws.Range("J22").Value = 999.8
$ws.Range("K22").Value = 1699.5
$ws.Range("L22").Value = 999.8
$ws.Range("M22").Value = -1404.5
$ws.Range("N22").Value = -1589.8
# Row 27
$ws.Range("H27").Value = 1199.7142
$ws.Range("I27").Value = 1699.5
$ws.Range("J27").Value = 999.8
$ws.Range("K27").Value = 1699.5
$ws.Range("L27").Value = 999.8
$ws.Range("M27").Value = -1592.5
$ws.Range("N27").Value = -1213.8
# Row 46
$ws.Range("H46").Value = 6089.5713
$ws.Range("I46").Value = 7102.4814
$ws.Range("J46").Value = 2671
$ws.Range("K46").Value = 7102.4814
$ws.Range("L46").Value = 2671
$ws.Range("M46").Value = -6914.4814
$ws.Range("N46").Value = -3047
# Row 132
$ws.Range("H132").Value = 24124.916
$ws.Range("I132").Value = 24160.1
$ws.Range("J132").Value = 23949
$ws.Range("K132").Value = 72480.29999999999
$ws.Range("L132").Value = 71847
$ws.Range("M132").Value = -69950.29999999999
$ws.Range("N132").Value = -76907

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 51
$ws.Range("H51").Value = 32000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 32000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 32000
$ws.Range("N51").Value = -33020
# Row 54
$ws.Range("H54").Value = 22800
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 22800
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 22800
$ws.Range("N54").Value = -23840
# Row 96
$ws.Range("H96").Value = 2896.5
$ws.Range("I96").Value = 3018.375
$ws.Range("J96").Value = 2815.25
$ws.Range("K96").Value = 3018.375
$ws.Range("L96").Value = 2815.25
$ws.Range("M96").Value = -1645.375
$ws.Range("N96").Value = -5561.25
# Row 122
$ws.Range("H122").Value = 3529.7297
$ws.Range("I122").Value = 3111.9656
$ws.Range("J122").Value = 5044.125
$ws.Range("K122").Value = 9335.8968
$ws.Range("L122").Value = 15132.375
$ws.Range("M122").Value = -6885.8968
$ws.Range("N122").Value = -20032.375
# Row 126
$ws.Range("H126").Value = 3540.9167
$ws.Range("I126").Value = 2499.182
$ws.Range("J126").Value = 15000
$ws.Range("K126").Value = 7497.545999999999
$ws.Range("L126").Value = 45000
$ws.Range("M126").Value = -5027.545999999999
$ws.Range("N126").Value = -49940
# Row 132
$ws.Range("H132").Value = 2794.5508
$ws.Range("I132").Value = 2701.9812
$ws.Range("J132").Value = 3101.1875
$ws.Range("K132").Value = 8105.943600000001
$ws.Range("L132").Value = 9303.5625
$ws.Range("M132").Value = -5575.943600000001
$ws.Range("N132").Value = -14363.5625

